# Refresh the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped figures. Cells whose new text would otherwise be
# auto-parsed by Excel as a number (and so lose e.g. trailing zeros) are
# forced to Text format first so the literal string is preserved verbatim.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.231.66"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.861.64"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7114"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.85"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08190"
$ws.Range("E8").Value = "  +10.80%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08167"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.888.65"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.173"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7090"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.68"
$ws.Range("D16").Value = "29.255.12"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007944"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.790"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.32"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D22").Value = "2.104.58"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.414"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.64"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1458"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.959"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.426"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.393"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7084"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01858"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").Value = "1.143.19"
$ws.Range("E41").Value = "  +6.76%  "
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.898"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.42"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.85"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.776"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").Value = "2.006.17"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.213"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.955"
$ws.Range("E51").Value = "  -1.27%  "
